$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (GitHub Actions scheduled update)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.671.29'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.953.03'
$ws.Range('E3').Value = '  +2.18%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '589.65'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.24'
$ws.Range('E6').Value = '  +4.18%  '
$ws.Range('E7').Value = '  -0.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.949.81'
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.505'
$ws.Range('E9').Value = '  +2.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.94'
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.149'
$ws.Range('E11').Value = '  +8.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.436'
$ws.Range('E12').Value = '  +1.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000232'
$ws.Range('E13').Value = '  +6.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.12'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.434.78'
$ws.Range('E16').Value = '  +1.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.589.67'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.947.70'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.64'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '435.20'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.43'
$ws.Range('E21').Value = '  +0.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.660'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.94'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.20'
$ws.Range('E24').Value = '  +5.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.13'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.89'
$ws.Range('E26').Value = '  +4.12%  '
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('E29').Value = '  +6.28%  '
$ws.Range('E30').Value = '  +2.93%  '
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0000100'
$ws.Range('E32').Value = '  +17.02%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.108'
$ws.Range('E33').Value = '  +2.20%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.14'
$ws.Range('E34').Value = '  -0.13%  '
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.989'
$ws.Range('E36').Value = '  +0.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.55'
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.00'
$ws.Range('E38').Value = '  +5.64%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.66'
$ws.Range('E39').Value = '  +0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.00'
$ws.Range('E40').Value = '  +3.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.35'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.115'
$ws.Range('E42').Value = '  -3.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.275'
$ws.Range('E43').Value = '  +3.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.40'
$ws.Range('E44').Value = '  -3.29%  '
$ws.Range('B45').Value = 'Monero'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '134.63'
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.678.61'
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0335'
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '352.55'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.104'
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.48'
$ws.Range('E51').Value = '  -1.57%  '
